$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the "mean" and "ptp" columns (D and E), including header labels,
# for rows 1 (header) through 11 (last data row).
$tmp = $ws.Range("D1:D11").Value2
$ws.Range("D1:D11").Value2 = $ws.Range("E1:E11").Value2
$ws.Range("E1:E11").Value2 = $tmp
